# Add 2022-Q4 data:
#  1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. before "2022-Q2"),
#     and fill it with the new quarter's fund-holding detail rows.
#  2. Update the "总计" summary sheet with a new row for 2022-Q4 (shifting the
#     existing quarterly rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet positioned right after "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet    = $wb.Worksheets.Item("2022-Q2")

$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# Helper: write a text-typed cell (force text storage, matching the workbook's
# existing convention where fund code / size / weight columns are stored as
# text even though they look numeric).
function Set-TextCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = [string]$val
}

function Set-NumberCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# Header row (row 1) - same headers/layout as the other quarterly sheets.
Set-TextCell $newSheet 1 2 "基金代码"
Set-TextCell $newSheet 1 3 "基金名称"
Set-TextCell $newSheet 1 4 "基金规模"
Set-TextCell $newSheet 1 5 "股票总仓位"
Set-TextCell $newSheet 1 6 "仓位占比"
Set-TextCell $newSheet 1 7 "持有市值(亿元)"
Set-TextCell $newSheet 1 8 "仓位排名"

# Data rows (row 2..12) — columns: A=index(n) B=code(t) C=name(t) D=size(t)
# E=stock position(t) F=position ratio(t) G=held value(t) H=rank(n)
$rows = @(
    @{ A=0;  B="513060"; C="博时恒生医疗保健ETF（QDII）";            D="69.51"; E="99.65"; F="5.17"; G="3.5937"; H=4 },
    @{ A=1;  B="159892"; C="华夏恒生香港上市生物科技ETF（QDII）";      D="5.36";  E="99.47"; F="5.70"; G="0.3055"; H=4 },
    @{ A=2;  B="513700"; C="鹏华中证港股通医药卫生综合ETF";           D="4.34";  E="95.37"; F="3.99"; G="0.1732"; H=5 },
    @{ A=3;  B="513120"; C="广发中证香港创新药（QDII-ETF）";          D="1.85";  E="98.69"; F="7.66"; G="0.1417"; H=4 },
    @{ A=4;  B="517380"; C="天弘恒生沪深港创新药精选50ETF";           D="2.05";  E="99.22"; F="3.98"; G="0.0816"; H=6 },
    @{ A=5;  B="501021"; C="华宝标普香港上市中国中小盘指数（LOF）A";   D="4.74";  E="94.57"; F="1.66"; G="0.0787"; H=7 },
    @{ A=6;  B="513280"; C="汇添富恒生香港上市生物科技ETF（QDII）";    D="1.35";  E="94.55"; F="5.48"; G="0.0740"; H=4 },
    @{ A=7;  B="513200"; C="易方达中证港股通医药卫生综合ETF";          D="1.69";  E="95.67"; F="4.11"; G="0.0695"; H=5 },
    @{ A=8;  B="159776"; C="银华中证港股通医药卫生综合ETF";            D="0.81";  E="93.98"; F="3.93"; G="0.0318"; H=5 },
    @{ A=9;  B="159718"; C="平安中证港股通医药卫生综合ETF";            D="0.67";  E="94.63"; F="3.94"; G="0.0264"; H=5 },
    @{ A=10; B="006127"; C="华宝标普香港上市中国中小盘指数（LOF）C";   D="0.45";  E="94.57"; F="1.66"; G="0.0075"; H=7 }
)

$r = 2
foreach ($row in $rows) {
    Set-NumberCell $newSheet $r 1 $row.A
    Set-TextCell   $newSheet $r 2 $row.B
    Set-TextCell   $newSheet $r 3 $row.C
    Set-TextCell   $newSheet $r 4 $row.D
    Set-TextCell   $newSheet $r 5 $row.E
    Set-TextCell   $newSheet $r 6 $row.F
    Set-TextCell   $newSheet $r 7 $row.G
    Set-NumberCell $newSheet $r 8 $row.H
    $r++
}

# ---------------------------------------------------------------------------
# 2. Update "总计" (summary) sheet: push existing quarter rows down by one
#    and insert the new 2022-Q4 figures at the top (row 2).
# ---------------------------------------------------------------------------
$ts = $totalSheet

# Capture current rows 2..8 (oldest at bottom) before overwriting, reading
# from the bottom up so we never clobber a row we still need to read.
for ($row = 8; $row -ge 2; $row--) {
    $dateVal  = $ts.Cells.Item($row, 2).Value2
    $countVal = $ts.Cells.Item($row, 3).Value2
    $mktVal   = $ts.Cells.Item($row, 4).Value2

    $ts.Cells.Item($row + 1, 1).Value = $row - 1
    $ts.Cells.Item($row + 1, 2).Value = $dateVal
    $ts.Cells.Item($row + 1, 3).Value = $countVal
    $ts.Cells.Item($row + 1, 4).Value = $mktVal
}

# New first data row: 2022-Q4
$ts.Cells.Item(2, 1).Value = 0
$ts.Cells.Item(2, 2).Value = "2022-Q4"
$ts.Cells.Item(2, 3).Value = 11
$ts.Cells.Item(2, 4).Value = 4.58
